$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"
$ws.Range("D2").Value = 0.026
$ws.Range("G2").Value = 0.07937853107344633
$ws.Range("H2").Value = 0.07937853107344633
$ws.Range("I2").Value = -0.1926553672316384
$ws.Range("J2").Value = -0.1926553672316384
$ws.Range("K2").Value = -6.91
$ws.Range("L2").Value = -0.1951977401129944
$ws.Range("U2").Value = 5.35
$ws.Range("V2").Value = 0.03331257783312578
$ws.Range("W2").Value = -0.06066725197541703
$ws.Range("X2").Value = 0.0671426390441774
$ws.Range("Y2").Value = -0.1278098910195944
$ws.Range("Z2").Value = 0.1687562568527435
$ws.Range("AA2").Value = -0.03251179863660199
$ws.Range("AB2").Value = 0.05533213187662268
$ws.Range("AC2").Value = -0.08784393051322467
$ws.Range("AD2").Value = 121.2
$ws.Range("AF2").Value = 121.2
$ws.Range("AG2").Value = 115.85
$ws.Range("AH2").Value = 0.4300922640170334
$ws.Range("AI2").Value = 0.4751078008624069
$ws.Range("AJ2").Value = 0.4190631217218304
$ws.Range("AK2").Value = 0.4638638638638639
$ws.Range("AL2").Value = 4.74
$ws.Range("AM2").Value = 3.46
$ws.Range("AN2").Value = -41.93771626297578
$ws.Range("AO2").Value = -1.438818565400844
$ws.Range("AP2").Value = -40.08650519031142
$ws.Range("AQ2").Value = -1.971098265895954

# Row 3 updates
$ws.Range("D3").Value = 0.026
$ws.Range("G3").Value = 0.07937853107344633
$ws.Range("H3").Value = 0.07937853107344633
$ws.Range("I3").Value = -0.1926553672316384
$ws.Range("J3").Value = -0.1926553672316384
$ws.Range("K3").Value = -6.91
$ws.Range("L3").Value = -0.1951977401129944
$ws.Range("O3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("U3").Value = 5.35
$ws.Range("V3").Value = 0.03331257783312578
$ws.Range("W3").Value = -0.06066725197541703
$ws.Range("X3").Value = 0.0671426390441774
$ws.Range("Y3").Value = -0.1278098910195944
$ws.Range("Z3").Value = 0.1687562568527435
$ws.Range("AA3").Value = -0.03251179863660199
$ws.Range("AB3").Value = 0.05533213187662268
$ws.Range("AC3").Value = -0.08784393051322467
$ws.Range("AD3").Value = 121.2
$ws.Range("AF3").Value = 121.2
$ws.Range("AG3").Value = 115.85
$ws.Range("AH3").Value = 0.4300922640170334
$ws.Range("AI3").Value = 0.4751078008624069
$ws.Range("AJ3").Value = 0.4190631217218304
$ws.Range("AK3").Value = 0.4638638638638639
$ws.Range("AL3").Value = 4.74
$ws.Range("AM3").Value = 3.46
$ws.Range("AN3").Value = -41.93771626297578
$ws.Range("AO3").Value = -1.438818565400844
$ws.Range("AP3").Value = -40.08650519031142
$ws.Range("AQ3").Value = -1.971098265895954

# Remove row 4 entirely (Capital Concept Limited JSC entry removed)
$ws.Rows("4").Delete()
